$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells A2:B4 hold plain text values (not numbers), and new rows 5:6 are
# appended, also as text. Force text number-format first so numeric-looking
# strings ("2", "1", ...) are not auto-converted to numbers, then restore the
# default "Normal" style so no stray style index is left on the cells.
$dataRange = $ws.Range("A2:B6")
$dataRange.NumberFormat = "@"

$ws.Range("A2").Value = "2"
$ws.Range("B2").Value = "1"

$ws.Range("A3").Value = "3"
$ws.Range("B3").Value = "4"

$ws.Range("A4").Value = "4"
$ws.Range("B4").Value = "3"

$ws.Range("A5").Value = "5"
$ws.Range("B5").Value = "8"

$ws.Range("A6").Value = "הדס"
$ws.Range("B6").Value = "9"

$dataRange.Style = "Normal"
